$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (Volume number and report date range) ---
$a8 = $ws.Range("A8")
$a8.Characters(21, 2).Text = "29"

$c9 = $ws.Range("C9")
$c9.Characters(27, 8).Text = "7/15/2024"
$c9.Characters(47, 9).Text = "7/21/2024"

# --- Crime statistics table updates (rows 15-30) ---
$ws.Range("C15").Value = 1
$ws.Range("C15").NumberFormat = "#,##0"
$ws.Range("D15").Value = "0"
$ws.Range("D15").NumberFormat = "General"
$ws.Range("E15").Value = "***.*"
$ws.Range("E15").NumberFormat = "General"
$ws.Range("F15").Value = 3
$ws.Range("G15").Value = 2
$ws.Range("H15").Value = 50
$ws.Range("I15").Value = 12
$ws.Range("J15").Value = 9
$ws.Range("K15").Value = 33.333333333333
$ws.Range("L15").Value = 9.090909090909
$ws.Range("M15").Value = 9.090909090909
$ws.Range("N15").Value = -20
$ws.Range("C16").Value = 5
$ws.Range("D16").Value = 5
$ws.Range("E16").Value = 0
$ws.Range("F16").Value = 17
$ws.Range("G16").Value = 15
$ws.Range("H16").Value = 13.333333333333
$ws.Range("I16").Value = 118
$ws.Range("J16").Value = 106
$ws.Range("K16").Value = 11.320754716981
$ws.Range("L16").Value = 10.280373831775
$ws.Range("M16").Value = -7.086614173228
$ws.Range("N16").Value = -53.174603174603
$ws.Range("C17").Value = 3
$ws.Range("D17").Value = 10
$ws.Range("E17").Value = -70
$ws.Range("F17").Value = 21
$ws.Range("G17").Value = 34
$ws.Range("H17").Value = -38.235294117647
$ws.Range("I17").Value = 184
$ws.Range("J17").Value = 174
$ws.Range("K17").Value = 5.747126436781
$ws.Range("L17").Value = 30.496453900709
$ws.Range("M17").Value = 70.37037037037
$ws.Range("N17").Value = 33.333333333333
$ws.Range("C18").Value = 1
$ws.Range("D18").Value = 2
$ws.Range("E18").Value = -50
$ws.Range("F18").Value = 14
$ws.Range("G18").Value = 5
$ws.Range("H18").Value = 180
$ws.Range("I18").Value = 83
$ws.Range("J18").Value = 81
$ws.Range("K18").Value = 2.469135802469
$ws.Range("L18").Value = 43.103448275862
$ws.Range("M18").Value = -54.395604395604
$ws.Range("N18").Value = -83.164300202839
$ws.Range("C19").Value = 14
$ws.Range("D19").Value = 5
$ws.Range("E19").Value = 180
$ws.Range("F19").Value = 63
$ws.Range("G19").Value = 43
$ws.Range("H19").Value = 46.511627906976
$ws.Range("I19").Value = 372
$ws.Range("J19").Value = 370
$ws.Range("K19").Value = 0.54054054054
$ws.Range("L19").Value = 24
$ws.Range("M19").Value = 51.219512195122
$ws.Range("N19").Value = 62.445414847161
$ws.Range("C20").Value = 16
$ws.Range("D20").Value = 14
$ws.Range("E20").Value = 14.285714285714
$ws.Range("F20").Value = 35
$ws.Range("G20").Value = 43
$ws.Range("H20").Value = -18.60465116279
$ws.Range("I20").Value = 221
$ws.Range("J20").Value = 315
$ws.Range("K20").Value = -29.841269841269
$ws.Range("L20").Value = 46.357615894039
$ws.Range("M20").Value = 123.232323232323
$ws.Range("N20").Value = -80.143755615453
$ws.Range("C21").Value = 40
$ws.Range("D21").Value = 36
$ws.Range("E21").Value = 11.111111111111
$ws.Range("F21").Value = 153
$ws.Range("G21").Value = 143
$ws.Range("H21").Value = 6.993006993006
$ws.Range("I21").Value = 992
$ws.Range("J21").Value = 1057
$ws.Range("K21").Value = -6.149479659413
$ws.Range("L21").Value = 28.664072632944
$ws.Range("M21").Value = 28.331177231565
$ws.Range("N21").Value = -55.852247441032
$ws.Range("L22").Value = -63.636363636363
$ws.Range("C23").Value = "0"
$ws.Range("C23").NumberFormat = "General"
$ws.Range("D23").Value = 1
$ws.Range("E23").Value = -100
$ws.Range("F23").Value = 5
$ws.Range("G23").Value = 8
$ws.Range("H23").Value = -37.5
$ws.Range("I23").Value = 32
$ws.Range("J23").Value = 37
$ws.Range("K23").Value = -13.513513513513
$ws.Range("L23").Value = 52.380952380952
$ws.Range("M23").Value = 39.130434782608
$ws.Range("N23").Value = "***.*"
$ws.Range("C24").Value = 21
$ws.Range("D24").Value = 29
$ws.Range("E24").Value = -27.586206896551
$ws.Range("F24").Value = 92
$ws.Range("G24").Value = 114
$ws.Range("H24").Value = -19.298245614035
$ws.Range("I24").Value = 742
$ws.Range("J24").Value = 785
$ws.Range("K24").Value = -5.477707006369
$ws.Range("L24").Value = 17.591125198098
$ws.Range("M24").Value = -5.59796437659
$ws.Range("N24").Value = "***.*"
$ws.Range("C25").Value = 14
$ws.Range("D25").Value = 13
$ws.Range("E25").Value = 7.692307692307
$ws.Range("F25").Value = 52
$ws.Range("G25").Value = 64
$ws.Range("H25").Value = -18.75
$ws.Range("I25").Value = 419
$ws.Range("J25").Value = 453
$ws.Range("K25").Value = -7.505518763796
$ws.Range("L25").Value = 30.124223602484
$ws.Range("M25").Value = "***.*"
$ws.Range("N25").Value = "***.*"
$ws.Range("C26").Value = 11
$ws.Range("D26").Value = 11
$ws.Range("E26").Value = 0
$ws.Range("F26").Value = 48
$ws.Range("G26").Value = 42
$ws.Range("H26").Value = 14.285714285714
$ws.Range("I26").Value = 305
$ws.Range("J26").Value = 283
$ws.Range("K26").Value = 7.773851590106
$ws.Range("L26").Value = 17.760617760617
$ws.Range("M26").Value = 27.615062761506
$ws.Range("N26").Value = "***.*"
$ws.Range("C27").Value = 1
$ws.Range("C27").NumberFormat = "#,##0"
$ws.Range("D27").Value = 1
$ws.Range("E27").Value = 0
$ws.Range("F27").Value = 3
$ws.Range("G27").Value = 4
$ws.Range("H27").Value = -25
$ws.Range("I27").Value = 17
$ws.Range("J27").Value = 18
$ws.Range("K27").Value = -5.555555555555
$ws.Range("L27").Value = -19.047619047619
$ws.Range("M27").Value = "***.*"
$ws.Range("N27").Value = "***.*"
$ws.Range("C28").Value = "0"
$ws.Range("C28").NumberFormat = "General"
$ws.Range("D28").Value = 1
$ws.Range("D28").NumberFormat = "#,##0"
$ws.Range("E28").Value = -100
$ws.Range("E28").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("F28").Value = 4
$ws.Range("G28").Value = 1
$ws.Range("H28").Value = 300
$ws.Range("I28").Value = 40
$ws.Range("J28").Value = 22
$ws.Range("K28").Value = 81.818181818181
$ws.Range("L28").Value = 29.032258064516
$ws.Range("M28").Value = "***.*"
$ws.Range("N28").Value = "***.*"
$ws.Range("G29").Value = 4
$ws.Range("H29").Value = -75
$ws.Range("G30").Value = 3
$ws.Range("H30").Value = -66.666666666666
